$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 722.5
$ws.Range("I19").Value = 807.5
$ws.Range("J19").Value = 609.1667
$ws.Range("K19").Value = 807.5
$ws.Range("L19").Value = 609.1667
$ws.Range("M19").Value = -632.5
$ws.Range("N19").Value = -959.1667
$ws.Range("H29").Value = 3265.2144
$ws.Range("I29").Value = 905
$ws.Range("J29").Value = 4576.4443
$ws.Range("K29").Value = 2715
$ws.Range("L29").Value = 13729.3329
$ws.Range("M29").Value = -2434
$ws.Range("N29").Value = -14291.3329
$ws.Range("H86").Value = 2154.2
$ws.Range("I86").Value = 1856.4286
$ws.Range("J86").Value = 2414.75
$ws.Range("K86").Value = 1856.4286
$ws.Range("L86").Value = 2414.75
$ws.Range("M86").Value = -733.4286
$ws.Range("N86").Value = -4660.75
$ws.Range("H89").Value = 2154.2
$ws.Range("I89").Value = 1856.4286
$ws.Range("J89").Value = 2414.75
$ws.Range("K89").Value = 9282.143
$ws.Range("L89").Value = 12073.75
$ws.Range("M89").Value = -3666.143
$ws.Range("N89").Value = -23305.75
$ws.Range("H100").Value = 3655.923
$ws.Range("I100").Value = 2800
$ws.Range("K100").Value = 2800
$ws.Range("M100").Value = -2259
$ws.Range("H107").Value = 1608.1818
$ws.Range("I107").Value = 1601.6923
$ws.Range("K107").Value = 1601.6923
$ws.Range("M107").Value = 318.3077000000001
$ws.Range("H111").Value = 621.5
$ws.Range("J111").Value = 900
$ws.Range("L111").Value = 2700
$ws.Range("N111").Value = -8834
$ws.Range("H137").Value = 2925
$ws.Range("I137").Value = 2050.5
$ws.Range("K137").Value = 6151.5
$ws.Range("M137").Value = -3601.5
$ws.Range("H141").Value = 5845.923
$ws.Range("I141").Value = 1499
$ws.Range("J141").Value = 6636.273
$ws.Range("K141").Value = 4497
$ws.Range("L141").Value = 19908.819
$ws.Range("M141").Value = 683
$ws.Range("N141").Value = -30268.819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5048.7144
$ws.Range("I32").Value = 4190.205
$ws.Range("K32").Value = 4190.205
$ws.Range("M32").Value = -3903.205
$ws.Range("H45").Value = 50003804
$ws.Range("I45").Value = 100001064
$ws.Range("J45").Value = 6551.1
$ws.Range("K45").Value = 100001064
$ws.Range("L45").Value = 6551.1
$ws.Range("M45").Value = -100000687
$ws.Range("N45").Value = -7305.1
$ws.Range("H61").Value = 9833
$ws.Range("I61").Value = 6739.6
$ws.Range("K61").Value = 6739.6
$ws.Range("M61").Value = -6527.6
$ws.Range("H74").Value = 33336866
$ws.Range("J74").Value = 1998.75
$ws.Range("L74").Value = 1998.75
$ws.Range("N74").Value = -3746.75
$ws.Range("H77").Value = 33336866
$ws.Range("J77").Value = 1998.75
$ws.Range("L77").Value = 9993.75
$ws.Range("N77").Value = -18729.75
$ws.Range("H110").Value = 3074.25
$ws.Range("I110").Value = 1820.7778
$ws.Range("J110").Value = 6834.6665
$ws.Range("K110").Value = 1820.7778
$ws.Range("L110").Value = 6834.6665
$ws.Range("M110").Value = 224.2221999999999
$ws.Range("N110").Value = -10924.6665
$ws.Range("H132").Value = 3464.158
$ws.Range("I132").Value = 2845.4614
$ws.Range("K132").Value = 8536.3842
$ws.Range("M132").Value = -6006.3842
$ws.Range("H133").Value = 72630.5
$ws.Range("J133").Value = 72630.5
$ws.Range("L133").Value = 72630.5
$ws.Range("N133").Value = -77690.5
$ws.Range("H136").Value = 9833
$ws.Range("I136").Value = 6739.6
$ws.Range("K136").Value = 20218.8
$ws.Range("M136").Value = -17668.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5210.3335
$ws.Range("I86").Value = 3018.0908
$ws.Range("K86").Value = 3018.0908
$ws.Range("M86").Value = -1895.0908
$ws.Range("H89").Value = 5210.3335
$ws.Range("I89").Value = 3018.0908
$ws.Range("K89").Value = 15090.454
$ws.Range("M89").Value = -9474.454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2158.75
$ws.Range("J16").Value = 6006.5
$ws.Range("L16").Value = 6006.5
$ws.Range("N16").Value = -6580.5
$ws.Range("H113").Value = 2158.75
$ws.Range("J113").Value = 6006.5
$ws.Range("L113").Value = 6006.5
$ws.Range("N113").Value = -10346.5
$ws.Range("H133").Value = 63332.332
$ws.Range("I133").Value = 64997
$ws.Range("J133").Value = 62500
$ws.Range("K133").Value = 64997
$ws.Range("L133").Value = 62500
$ws.Range("M133").Value = -62467
$ws.Range("N133").Value = -67560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 401.33334
$ws.Range("J24").Value = 1000
$ws.Range("L24").Value = 3000
$ws.Range("N24").Value = -3460
$ws.Range("H131").Value = 5557460.5
$ws.Range("I131").Value = 2247.2222
$ws.Range("J131").Value = 13890280
$ws.Range("K131").Value = 6741.6666
$ws.Range("L131").Value = 41670840
$ws.Range("M131").Value = -1701.6666
$ws.Range("N131").Value = -41680920

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2976.1428
$ws.Range("I102").Value = 2986.8
$ws.Range("K102").Value = 2986.8
$ws.Range("M102").Value = -1364.8
$ws.Range("H107").Value = 1668.4286
$ws.Range("I107").Value = 851.2857
$ws.Range("K107").Value = 851.2857
$ws.Range("M107").Value = 1068.7143
$ws.Range("H113").Value = 4960.4
$ws.Range("I113").Value = 3307.3333
$ws.Range("J113").Value = 7440
$ws.Range("K113").Value = 3307.3333
$ws.Range("L113").Value = 7440
$ws.Range("M113").Value = -1137.3333
$ws.Range("N113").Value = -11780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6660.8
$ws.Range("I61").Value = 1426
$ws.Range("J61").Value = 10150.667
$ws.Range("K61").Value = 1426
$ws.Range("L61").Value = 10150.667
$ws.Range("M61").Value = -1224
$ws.Range("N61").Value = -10554.667
$ws.Range("H113").Value = 6660.8
$ws.Range("I113").Value = 1426
$ws.Range("J113").Value = 10150.667
$ws.Range("K113").Value = 1426
$ws.Range("L113").Value = 10150.667
$ws.Range("M113").Value = 744
$ws.Range("N113").Value = -14490.667
$ws.Range("H132").Value = 3400.8235
$ws.Range("I132").Value = 2088
$ws.Range("K132").Value = 6264
$ws.Range("M132").Value = -3734

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1218.5
$ws.Range("I100").Value = 453.7143
$ws.Range("K100").Value = 907.4286
$ws.Range("M100").Value = -366.4286
$ws.Range("H113").Value = 1064.875
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H126").Value = 3149.4
$ws.Range("I126").Value = 1926.25
$ws.Range("J126").Value = 3964.8333
$ws.Range("K126").Value = 5778.75
$ws.Range("L126").Value = 11894.4999
$ws.Range("M126").Value = -3308.75
$ws.Range("N126").Value = -16834.4999
